$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for rows 2-11, columns B-G (ME, MAE, MSE, RMSE, SE, N)
$data = @{
    2  = @(-0.0258731143856077, 0.8712523508600251, 1.357352004991913, 1.165054507305093, 1.196684520570948, 19)
    3  = @(-0.05889625796533703, 0.7809309034430243, 1.486890174830229, 1.219381062191073, 1.25326832255184, 18)
    4  = @(-0.03036361610786248, 0.7674964733547235, 0.9590082094163144, 0.9792896453125165, 1.008943334136136, 17)
    5  = @(0.09993964655644642, 0.6259185880536732, 0.8952003656645542, 0.9461502870393024, 0.9717132266730971, 16)
    6  = @(0.03831129326750842, 0.7119710808408777, 1.080030924484957, 1.039245362984583, 1.074989951428016, 15)
    7  = @(0.09938456961952043, 0.7378794059141224, 1.55852155670565, 1.248407608397854, 1.291421985981127, 14)
    8  = @(0.06981294092732569, 0.5769263726184313, 0.6329722743072022, 0.7955955469377655, 0.8248878488730877, 13)
    9  = @(0.2198501169052137, 0.8740014322694108, 1.535812380474684, 1.239278976048042, 1.273853840583311, 12)
    10 = @(0.1257302919151607, 0.7877651198835861, 0.9791626505880603, 0.9895264779620908, 1.029412453155336, 11)
    11 = @(0.2215020228704644, 0.8099270631829588, 1.114429208763017, 1.055665292014006, 1.08799819374202, 10)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
